$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 48..79 in the "Artfynd" sheet get cyclically rotated down by one:
#   new row 49 = old row 48
#   new row 50 = old row 49
#   ...
#   new row 79 = old row 78
#   new row 48 = old row 79   (wrap-around)
#
# We use Range.Copy (not .Value assignment) on the used column span
# (A:AY) so cell types are preserved verbatim (e.g. literal date-like
# text such as "2017-09-29" is not re-interpreted as a real
# date/number), and no stray style records get introduced. We avoid
# whole-row Range objects (Rows.Item) since copying full 16384-column
# rows is drastically slower than copying the actual used range.
#
# NOTE: this engine's Range.Copy behaves like "paste, skip blanks" -
# copying a blank source cell onto a non-blank destination leaves the
# destination untouched. So every destination row is explicitly
# Clear()-ed immediately before each Copy() to guarantee blanks really
# overwrite whatever was there.

$firstRow = 48
$lastRow = 79
$scratchRow = 81
$firstCol = "A"
$lastCol = "AY"

function RowRange($r) {
    return $ws.Range($firstCol + $r + ":" + $lastCol + $r)
}

# 1) Stash the last row (79) in a scratch row far below the used range.
(RowRange $scratchRow).Clear()
(RowRange $lastRow).Copy((RowRange $scratchRow))

# 2) Shift rows firstRow..lastRow-1 down by one, working from the bottom
#    up so we never clobber a row before we've read it as a source.
for ($r = $lastRow - 1; $r -ge $firstRow; $r--) {
    (RowRange ($r + 1)).Clear()
    (RowRange $r).Copy((RowRange ($r + 1)))
}

# 3) Drop the stashed old row 79 into row 48 (the wrap-around).
(RowRange $firstRow).Clear()
(RowRange $scratchRow).Copy((RowRange $firstRow))

# 4) Clean up the scratch row so it doesn't leave stray data behind.
(RowRange $scratchRow).Clear()

$excel.CutCopyMode = $false
